$wb = $excel.ActiveWorkbook

# --- Sheet: Power Consumption Est. (4th sheet) ---
$ws = $wb.Worksheets.Item("Power Consumption Est.")

# Header row
$ws.Range("C1").Value = "ON-Time(min)"
$ws.Range("D1").Value = "Energy Used (Wh)"

# Rows 2-7: component power table
$ws.Range("A2").Value = "HASP"
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = "Radio/Antenna"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "OBC"
$ws.Range("B4").Value = 0.0004

$ws.Range("A5").Value = "GPS"
$ws.Range("B5").Value = 0.5

$ws.Range("A6").Value = "Power Mang."
$ws.Range("B6").Value = 0.000115

$ws.Range("A7").Value = "Attitude"
$ws.Range("B7").Value = 0

# Clear old row 8 (Batteries) that is no longer there
$ws.Range("A8").ClearContents()

# Row 9: Total
$ws.Range("A9").Value = "Total "
$ws.Range("B9").Formula = "=SUM(B2:B7)"

# Row 11: Batteries
$ws.Range("A11").Value = "Batteries"
$ws.Range("B11").Value = 0

# Row 12: Solar 1 side
$ws.Range("A12").Value = "Solar 1 side"
$ws.Range("B12").Value = 6.62

# Row 14: Power in - out
$ws.Range("A14").Value = "Power in - out"
$ws.Range("B14").Formula = "=B12-(B9+B11)"

# E11: Total Battery wattage to be determined (set last so it becomes the last new shared string)
$ws.Range("E11").Value = "Total Battery wattage to be determined "

# Column widths (A: 12.42578125 -> 15.42578125)
$ws.Columns.Item(1).ColumnWidth = 15.42578125

# Selection / active cell
[void]$ws.Range("E16").Select()

# --- Sheet: Solar Panels (1st sheet) selection change ---
$ws1 = $wb.Worksheets.Item("Solar Panels")
[void]$ws1.Range("J2").Select()

# --- Sheet: Power Managment (3rd sheet) selection change ---
$ws3 = $wb.Worksheets.Item("Power Managment")
[void]$ws3.Range("H8").Select()

# Re-activate the Power Consumption Est. sheet and set final selection
[void]$ws.Activate()
[void]$ws.Range("E16").Select()
